$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: columns C/D/E are reordered ---
# old: C1=max, D1=prediction, E1=rejection-f
# new: C1=prediction, D1=rejection-f, E1=max
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# --- Data rows 2-9 ---
# old: C=max numeric value, D=species string, E=species string
# new: C=species string, D=species string (unchanged), E=new numeric "max" value
$species = "s__Bacteroides_F pectinophilus"

$newMax = @{
    2 = 0.6129018169920266
    3 = 0.6156466526235611
    4 = 0.6152407737837783
    5 = 0.6135368067503487
    6 = 0.6194687776471578
    7 = 0.6175214631405851
    8 = 0.6126562641101683
    9 = 0.61012922088622
}

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = $species
    $ws.Cells.Item($r, 4).Value = $species
    $ws.Cells.Item($r, 5).Value = $newMax[$r]
}
